# Edit script: apply the Notes.docx sprint-1 change.
#
# Summary of changes:
#  1. Remove the <w:bookmarkStart/><w:bookmarkEnd/> pair (id=0, name=_GoBack)
#     that currently sits right after the "Output:" heading run.
#  2. In the "Data" header cell of the Output table, add a right tab stop
#     at 4459 twips to the paragraph, and append a new (empty-of-visible-
#     text) run ending in a tab character with the same bold/size run
#     formatting as "Data".
#  3. Clean up several JSON-ish lines in the Output table's first cell:
#     merge runs that used to be split apart by now-removed
#     w:proofErr gramStart/gramEnd markers into single runs with the
#     same combined text (no visible/textual change, just simplifies
#     the run/proofErr structure) for the FromBay/ToBay/MaxRow/MaxTier/
#     IsVirtual/Href lines.
#  4. On the "Rel" line, change the value from "block-information" to
#     "cell-templates", split across 3 runs, and re-insert the
#     bookmarkStart/bookmarkEnd (_GoBack) pair between the 2nd and 3rd
#     of those runs (i.e. right after the new "cell-templates" text).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: remove the _GoBack bookmark after "Output:"
# ---------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------
# Step 2: add tab stop + trailing tab run to the "Data" header cell
# ---------------------------------------------------------------------
$outTable = $d.Tables(4)
$dataPara = $outTable.Cell(1, 1).Range.Paragraphs(1)
$dataPara.Range.ParagraphFormat.TabStops.Add(222.95, 2)

$dataCellRange = $outTable.Cell(1, 1).Range
$dataCellRange.Collapse(0)
$dataCellRange.MoveEnd(1, -1)
$dataCellRange.InsertAfter("`t")

# ---------------------------------------------------------------------
# Step 3: simplify the JSON lines (merge split runs / drop proofErr)
# ---------------------------------------------------------------------
$jsonCell = $outTable.Cell(2, 1).Range

$replacements = @(
    @('   “FromBay”:<frombay>,', '   “FromBay”:<frombay>,'),
    @('    “ToBay”:<tobay>,', '    “ToBay”:<tobay>,'),
    @('    “MaxRow”:<maxrow>,', '    “MaxRow”:<maxrow>,'),
    @('    “MaxTier”:<maxtier>,', '    “MaxTier”:<maxtier>,'),
    @('    “IsVirtual”:<isvirtual>,', '    “IsVirtual”:<isvirtual>,'),
    @('                    “Href”: ”', '                    “Href”: ”')
)

foreach ($pair in $replacements) {
    $jsonCell.Find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)
}

# ---------------------------------------------------------------------
# Step 4: change "block-information" -> "cell-templates" and split the
# "Rel" line into 3 runs with the _GoBack bookmark reinserted between
# the 2nd and 3rd.
# ---------------------------------------------------------------------
$jsonCell.Find.Execute('                    “Rel”: “block-information”,', $false, $false, $false, $false, $false, $true, 1, $false, '                    “Rel”: “cell-templates”,', 2)

$relFind = $jsonCell.Duplicate
$relFind.Find.Execute('cell-templates', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$relFind.Collapse(0)
$d.Bookmarks.Add("_GoBack", $relFind)
